$wb = $excel.ActiveWorkbook

# The same updates must be applied to both the "展览" sheet and the
# "全部类型" sheet, which both contain duplicated event data.
$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F ("想去人数")
$updates = @{
    7  = 1339
    8  = 465
    9  = 92
    11 = 115
    12 = 164
    15 = 131
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
